# resolve_report_task.xlsx — remove form_id from the "settings" sheet.
#
# The form_id column (column B on the "settings" sheet) is dropped, the
# remaining columns (version/style/namespaces) shift left, the cell
# comments documenting each column shift accordingly, and the
# conditionalFormatting ranges on the "survey" sheet get consolidated.

$wb = $excel.ActiveWorkbook
$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- 1. Shift the header comments on the settings sheet one column left ---
# B1 (form_id) is being removed; C1 (version) -> B1, D1 (style/pages) -> C1,
# E1 (namespaces) -> D1.
$versionText    = $settings.Range("C1").Comment.Text()
$pagesText      = $settings.Range("D1").Comment.Text()
$namespaceText  = $settings.Range("E1").Comment.Text()

$settings.Range("B1").Comment.Text($versionText)   | Out-Null
$settings.Range("C1").Comment.Text($pagesText)     | Out-Null
$settings.Range("D1").Comment.Text($namespaceText) | Out-Null
$settings.Range("E1").Comment.Delete()             | Out-Null

# --- 2. Delete the form_id column itself (column B) ---
$settings.Columns.Item(2).Delete() | Out-Null

# --- 3. Re-select B1 on the settings sheet without changing which tab is
#        the active/selected one (the "survey" sheet stays selected) ---
$settings.Range("B1").Select() | Out-Null
$survey.Activate() | Out-Null

# --- 4. Consolidate the conditionalFormatting ranges on the survey sheet ---
$rules = $survey.Cells.FormatConditions
for ($i = 1; $i -le $rules.Count; $i++) {
    $rule = $rules.Item($i)
    if ($rule.Formula1 -eq '=$A2="begin_group"') {
        $rule.ModifyAppliesToRange($survey.Range("A2:G10006")) | Out-Null
    }
    if ($rule.Formula1 -like '*ISBLANK(C2)*') {
        $rule.ModifyAppliesToRange($survey.Range("C2:C10006")) | Out-Null
    }
}
